$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A5").Value = "pedro"
$ws.Range("B5").Value = "donpedro"
[void]$ws.Range("A6").Select()
